$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")
$ws.Activate()

# --- New data values added to the comparison table ---

# Arabidopsis thaliana / EST row: number of promoters-TSRchitect
$ws.Range("H6").Value = 16238
$ws.Range("H6").NumberFormat = "#,##0"

# Zea mays / EST row: number of input reads
$ws.Range("D8").Value = 2019694

# Zea mays / cDNA row: number of input reads
$ws.Range("D9").Value = 27455

# Fill in alignment-file locations for the EST/cDNA rows (J9 first so the
# Zea mays string lands on shared-string index 42, matching A_thaliana at 43)
$ws.Range("J9").Value = "/projects/TSRplants/ESTcDNA/Z_mays/alignments/cDNAzm_?.fa.gsq; /projects/TSRplants/ESTcDNA/Z_mays/alignments/cDNAzm_?.bed"
$ws.Range("J8").Value = "/projects/TSRplants/ESTcDNA/A_thaliana/alignments/genbank_EST/AtESTgenbank_align.bed"

# --- View-state updates (scroll position / active selection) ---
$ws.Range("F9").Select() | Out-Null
